$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles) of the last existing data row (114) down to
# the two new rows (115, 116) so the new cells inherit the same per-column
# styles (bold/border on col A, custom date format on col E, default on the
# rest) without introducing any new style entries.
$ws.Range("A114:V114").Copy()
$ws.Range("A115:V116").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 115 ----
$ws.Cells.Item(115, 1).Value = 114
$ws.Cells.Item(115, 2).Value = "romania"
$ws.Cells.Item(115, 3).Value = "liga-1"
$ws.Cells.Item(115, 4).Value = "2023-2024"
$ws.Cells.Item(115, 5).Value = 45234.65625
$ws.Cells.Item(115, 6).Value = "Sepsi Sf. Gheorghe"
$ws.Cells.Item(115, 7).Value = 0
$ws.Cells.Item(115, 8).Value = "Petrolul"
$ws.Cells.Item(115, 9).Value = 0
$ws.Cells.Item(115, 10).Value = 1.75
$ws.Cells.Item(115, 11).Value = "02/11/2023 14:42"
$ws.Cells.Item(115, 12).Value = 1.89
$ws.Cells.Item(115, 13).Value = "04/11/2023 15:42"
$ws.Cells.Item(115, 14).Value = 3.49
$ws.Cells.Item(115, 15).Value = "02/11/2023 14:42"
$ws.Cells.Item(115, 16).Value = 3.26
$ws.Cells.Item(115, 17).Value = "04/11/2023 15:21"
$ws.Cells.Item(115, 18).Value = 4.6
$ws.Cells.Item(115, 19).Value = "02/11/2023 14:42"
$ws.Cells.Item(115, 20).Value = 4.74
$ws.Cells.Item(115, 21).Value = "04/11/2023 15:42"
$ws.Cells.Item(115, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/sepsi-petrolul/tbOPept6/"

# ---- Row 116 ----
$ws.Cells.Item(116, 1).Value = 115
$ws.Cells.Item(116, 2).Value = "romania"
$ws.Cells.Item(116, 3).Value = "liga-1"
$ws.Cells.Item(116, 4).Value = "2023-2024"
$ws.Cells.Item(116, 5).Value = 45234.84375
$ws.Cells.Item(116, 6).Value = "Univ. Craiova"
$ws.Cells.Item(116, 7).Value = 1
$ws.Cells.Item(116, 8).Value = "U Craiova 1948"
$ws.Cells.Item(116, 9).Value = 1
$ws.Cells.Item(116, 10).Value = 1.73
$ws.Cells.Item(116, 11).Value = "02/11/2023 15:28"
$ws.Cells.Item(116, 12).Value = 1.59
$ws.Cells.Item(116, 13).Value = "04/11/2023 20:12"
$ws.Cells.Item(116, 14).Value = 3.54
$ws.Cells.Item(116, 15).Value = "02/11/2023 15:28"
$ws.Cells.Item(116, 16).Value = 4.12
$ws.Cells.Item(116, 17).Value = "04/11/2023 20:12"
$ws.Cells.Item(116, 18).Value = 4.62
$ws.Cells.Item(116, 19).Value = "02/11/2023 15:28"
$ws.Cells.Item(116, 20).Value = 5.63
$ws.Cells.Item(116, 21).Value = "04/11/2023 20:05"
$ws.Cells.Item(116, 22).Value = "https://www.betexplorer.com/football/romania/liga-1/univ-craiova-fc-u-craiova/KrzHc6Bg/"
